$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("F2").Value = 1.51
$ws.Range("G2").Value = 1.65
$ws.Range("H2").Value = 6.8
$ws.Range("I2").Value = 10.5
$ws.Range("J2").Value = 3.65
$ws.Range("K2").Value = 5
$ws.Range("M2").Value = 1.08
$ws.Range("N2").Value = 3.2
$ws.Range("S2").Value = 3.35
$ws.Range("U2").Value = 1.7
$ws.Range("W2").Value = 2.52

# Row 3
$ws.Range("G3").Value = 2.02
$ws.Range("J3").Value = 3.2
$ws.Range("V3").Value = 1.22
$ws.Range("W3").Value = 1.98
$ws.Range("AN3").Value = 22

# Row 4
$ws.Range("G4").Value = 1.85
$ws.Range("J4").Value = 3.3
$ws.Range("K4").Value = 3.7
$ws.Range("O4").Value = 1.49
$ws.Range("Q4").Value = 2.42
$ws.Range("S4").Value = 4.4
$ws.Range("V4").Value = 1.17
$ws.Range("W4").Value = 2.16

# Row 5
$ws.Range("T5").Value = 1.97

# Row 6
$ws.Range("F6").Value = 1.24
$ws.Range("G6").Value = 1.25
$ws.Range("H6").Value = 14
$ws.Range("K6").Value = 7.8
$ws.Range("N6").Value = 8.199999999999999
$ws.Range("P6").Value = 3.4
$ws.Range("Q6").Value = 1.4
$ws.Range("R6").Value = 1.98
$ws.Range("S6").Value = 1.99
$ws.Range("T6").Value = 1.89
$ws.Range("U6").Value = 2.04
$ws.Range("W6").Value = 5
$ws.Range("Y6").Value = 65
$ws.Range("AA6").Value = 570
$ws.Range("AB6").Value = 14
$ws.Range("AE6").Value = 190
$ws.Range("AF6").Value = 9.6
$ws.Range("AN6").Value = 3.2
$ws.Range("AO6").Value = 160

# Row 7
$ws.Range("F7").Value = 3.35
$ws.Range("G7").Value = 3.4
$ws.Range("P7").Value = 2.28
$ws.Range("T7").Value = 1.63
$ws.Range("V7").Value = 1.75
$ws.Range("Y7").Value = 12.5
$ws.Range("AF7").Value = 24
$ws.Range("AO7").Value = 14

# Row 8
$ws.Range("J8").Value = 5
$ws.Range("Q8").Value = 1.82
$ws.Range("S8").Value = 3.1
$ws.Range("AB8").Value = 8
$ws.Range("AG8").Value = 9.800000000000001
$ws.Range("AO8").Value = 160

# Row 9
$ws.Range("F9").Value = 3.1
$ws.Range("G9").Value = 3.15
$ws.Range("H9").Value = 2.36
$ws.Range("I9").Value = 2.38
$ws.Range("J9").Value = 3.85
$ws.Range("K9").Value = 3.9
$ws.Range("P9").Value = 2.66
$ws.Range("U9").Value = 2.84
$ws.Range("V9").Value = 1.72
$ws.Range("AF9").Value = 25

# Row 10
$ws.Range("N10").Value = 6
$ws.Range("Q10").Value = 1.56
$ws.Range("T10").Value = 1.53
$ws.Range("U10").Value = 2.84
$ws.Range("X10").Value = 24
$ws.Range("AC10").Value = 9.6
$ws.Range("AG10").Value = 11
$ws.Range("AO10").Value = 18

# Row 11
$ws.Range("H11").Value = 3.6
$ws.Range("N11").Value = 4.9
$ws.Range("Q11").Value = 1.72
$ws.Range("T11").Value = 1.64
$ws.Range("W11").Value = 1.84
$ws.Range("AK11").Value = 19.5

# Row 12
$ws.Range("H12").Value = 19
$ws.Range("I12").Value = 19.5
$ws.Range("J12").Value = 9.800000000000001
$ws.Range("K12").Value = 10
$ws.Range("R12").Value = 2.46
$ws.Range("S12").Value = 1.65
$ws.Range("T12").Value = 1.79
$ws.Range("U12").Value = 2.22
$ws.Range("W12").Value = 6.4
$ws.Range("X12").Value = 100
$ws.Range("Y12").Value = 970
$ws.Range("Z12").Value = 250
$ws.Range("AD12").Value = 65
$ws.Range("AE12").Value = 260
$ws.Range("AF12").Value = 13
$ws.Range("AG12").Value = 14
$ws.Range("AH12").Value = 34
$ws.Range("AJ12").Value = 11.5
$ws.Range("AL12").Value = 29
$ws.Range("AN12").Value = 2.46

# Row 13
$ws.Range("I13").Value = 2.44
$ws.Range("L13").Value = 1.33
$ws.Range("O13").Value = 1.25
$ws.Range("P13").Value = 2.26
$ws.Range("Q13").Value = 1.76
$ws.Range("R13").Value = 1.51
$ws.Range("S13").Value = 2.86
$ws.Range("T13").Value = 1.63
$ws.Range("U13").Value = 2.5
$ws.Range("V13").Value = 1.69
$ws.Range("Y13").Value = 13
$ws.Range("AB13").Value = 15.5
$ws.Range("AE13").Value = 23
$ws.Range("AF13").Value = 23
$ws.Range("AM13").Value = 70
$ws.Range("AN13").Value = 24

# Row 14
$ws.Range("F14").Value = 2.24
$ws.Range("G14").Value = 2.54
$ws.Range("H14").Value = 2.84
$ws.Range("L14").Value = 1.28
$ws.Range("M14").Value = 1.05
$ws.Range("N14").Value = 4.5
$ws.Range("O14").Value = 1.24
$ws.Range("P14").Value = 2.26
$ws.Range("Q14").Value = 1.64
$ws.Range("R14").Value = 1.52
$ws.Range("S14").Value = 2.56
$ws.Range("T14").Value = 1.56
$ws.Range("U14").Value = 2.36
$ws.Range("V14").Value = 1.43
$ws.Range("W14").Value = 1.64
$ws.Range("X14").Value = 23
$ws.Range("Y14").Value = 16.5
$ws.Range("Z14").Value = 28
$ws.Range("AA14").Value = 60
$ws.Range("AB14").Value = 14.5
$ws.Range("AC14").Value = 10.5
$ws.Range("AD14").Value = 14
$ws.Range("AE14").Value = 34
$ws.Range("AF14").Value = 18.5
$ws.Range("AG14").Value = 12.5
$ws.Range("AH14").Value = 17
$ws.Range("AK14").Value = 1000
$ws.Range("AN14").Value = 17.5
$ws.Range("AO14").Value = 27
